$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# --- New product rows 2-11 (Model, EngineTypeId, CarTypeId, Horse Power, ManufacturerId, Release Year, Price) ---
$ws.Range("B2").Value = "BMW E93"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 232
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2011
$ws.Range("H2").Value = 56000

$ws.Range("B3").Value = "Honda Integra Type R"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 156
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 1998
$ws.Range("H3").Value = 36000

$ws.Range("B4").Value = "BMW X5"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 188
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1999
$ws.Range("H4").Value = 44000

$ws.Range("B5").Value = "Mercedes-Benz CLK"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 160
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 2008
$ws.Range("H5").Value = 68000

$ws.Range("B6").Value = "Toyota GT86"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 220
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 2012
$ws.Range("H6").Value = 130000

$ws.Range("B7").Value = "BMW 5-Series E39"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 145
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1995
$ws.Range("H7").Value = 32000

$ws.Range("B8").Value = "Audi A2"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 89
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 1999
$ws.Range("H8").Value = 29000

$ws.Range("B9").Value = "Toyota Hilux"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 166
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 1997
$ws.Range("H9").Value = 52000

$ws.Range("B10").Value = "Nissan GT-R"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 356
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 2008
$ws.Range("H10").Value = 110000

$ws.Range("B11").Value = "Toyota Prius"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 68
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 2004
$ws.Range("H11").Value = 70000

# Price column gets a thousands-separator number format
$ws.Range("H2:H11").NumberFormat = "#,##0"

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 26.29   # col B widened for longer model names
$ws.Columns.Item(8).ColumnWidth = 8.92    # new col H (Price) width

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Sheet selections ---
$dealers = $wb.Worksheets.Item("Dealers")
$dealers.Range("B11").Select()

$ws.Activate()
$ws.Range("C10").Select()
